$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:18:07'
$ws.Cells.Item(3, 1).Value = 'Total filas: 75'

$rows = @(
    @(6, '04:03:00', '04:03', '81_EL PELIGRO', 0, 'LP1912'),
    @(7, '04:37:19', '04:46', '215A_EL PATO', 9, 'LP1912'),
    @(8, '04:52:25', '04:53', '11_ETCHEVERRY', 1, 'LP1912'),
    @(9, '04:52:25', '05:16', '17_ROMERO', 24, 'LP1912'),
    @(10, '05:20:30', '05:20', '17_ROMERO', 0, 'LP1912'),
    @(11, '04:52:25', '05:22', '23_HERNANDEZ', 30, 'LP1912'),
    @(12, '05:20:30', '05:26', '23_HERNANDEZ', 6, 'LP1912'),
    @(13, '05:20:30', '05:34', '215B_EL PATO', 14, 'LP1912'),
    @(14, '04:03:00', '05:35', '215B_EL PATO', 92, 'LP1912'),
    @(15, '04:03:00', '05:41', '14_ABASTO', 98, 'LP1912'),
    @(16, '05:20:30', '05:46', '15_ABASTO', 26, 'LP1912'),
    @(17, '05:54:55', '05:54', '10_OLMOS', 0, 'LP1912'),
    @(18, '05:54:55', '05:55', '81_EL PELIGRO', 1, 'LP1912'),
    @(19, '05:20:30', '06:04', '16_SANTA ANA', 44, 'LP1912'),
    @(20, '05:54:55', '06:11', '215A_EL PATO', 17, 'LP1912'),
    @(21, '05:54:55', '06:13', '225_HARAS DEL SUR', 19, 'LP1912'),
    @(22, '05:20:30', '06:14', '225_HARAS DEL SUR', 54, 'LP1912'),
    @(23, '05:54:55', '06:20', '26_HERNANDEZ', 26, 'LP1912'),
    @(24, '05:20:30', '06:21', '26_HERNANDEZ', 61, 'LP1912'),
    @(25, '05:54:55', '06:26', '23_HERNANDEZ', 32, 'LP1912'),
    @(26, '06:24:16', '06:27', '23_HERNANDEZ', 3, 'LP1912'),
    @(27, '05:54:55', '06:29', '86_EST CHICA-ESC AGRARIA', 35, 'LP1912'),
    @(28, '06:24:16', '06:30', '86_EST CHICA-ESC AGRARIA', 6, 'LP1912'),
    @(29, '06:24:16', '06:31', '16_SANTA ANA', 7, 'LP1912'),
    @(30, '05:54:55', '06:43', '225_C ROCA-H SUR', 49, 'LP1912'),
    @(31, '06:24:16', '06:44', '225_C ROCA-H SUR', 20, 'LP1912'),
    @(32, '06:24:16', '06:46', '215C_EL PATO', 22, 'LP1912'),
    @(33, '06:53:31', '06:56', '14_ABASTO', 3, 'LP1912'),
    @(34, '05:54:55', '06:59', '14_ABASTO', 65, 'LP1912'),
    @(35, '06:24:16', '07:00', '14_ABASTO', 36, 'LP1912'),
    @(36, '06:53:31', '07:01', '16_SANTA ANA', 8, 'LP1912'),
    @(37, '06:53:31', '07:04', '23_HERNANDEZ', 11, 'LP1912'),
    @(38, '06:53:31', '07:05', '15_ABASTO', 12, 'LP1912'),
    @(39, '05:54:55', '07:06', '225_GOMEZ', 72, 'LP1912'),
    @(40, '06:53:31', '07:07', '225_GOMEZ', 14, 'LP1912'),
    @(41, '06:53:31', '07:11', '215A_EL PATO', 18, 'LP1912'),
    @(42, '06:53:31', '07:15', '11_ETCHEVERRY', 22, 'LP1912'),
    @(43, '06:53:31', '07:16', '16_SANTA ANA', 23, 'LP1912'),
    @(44, '06:24:16', '07:16', '11_ETCHEVERRY', 52, 'LP1912'),
    @(45, '07:18:07', '07:18', '16_SANTA ANA', 0, 'LP1912'),
    @(46, '05:54:55', '07:20', '26_HERNANDEZ', 86, 'LP1912'),
    @(47, '07:18:07', '07:21', '225_GOMEZ', 3, 'LP1912'),
    @(48, '07:18:07', '07:21', '10_OLMOS', 3, 'LP1912'),
    @(49, '07:18:07', '07:21', '26_HERNANDEZ', 3, 'LP1912'),
    @(50, '05:54:55', '07:22', '10_OLMOS', 88, 'LP1912'),
    @(51, '06:53:31', '07:23', '10_OLMOS', 30, 'LP1912'),
    @(52, '06:53:31', '07:31', '11_ETCHEVERRY', 38, 'LP1912'),
    @(53, '05:54:55', '07:31', '16_SANTA ANA', 97, 'LP1912'),
    @(54, '07:18:07', '07:32', '84_COLONIA URQUIZA-ESC 49', 14, 'LP1912'),
    @(55, '07:18:07', '07:32', '11_ETCHEVERRY', 14, 'LP1912'),
    @(56, '07:18:07', '07:32', '16_SANTA ANA', 14, 'LP1912'),
    @(57, '06:53:31', '07:36', '27_EL RETIRO', 43, 'LP1912'),
    @(58, '07:18:07', '07:36', '23_HERNANDEZ', 18, 'LP1912'),
    @(59, '07:18:07', '07:37', '27_EL RETIRO', 19, 'LP1912'),
    @(60, '06:24:16', '07:39', '10_OLMOS', 75, 'LP1912'),
    @(61, '06:53:31', '07:47', '14_ABASTO', 54, 'LP1912'),
    @(62, '07:18:07', '07:48', '14_ABASTO', 30, 'LP1912'),
    @(63, '06:53:31', '07:51', '215D_EL PATO', 58, 'LP1912'),
    @(64, '07:18:07', '07:52', '215D_EL PATO', 34, 'LP1912'),
    @(65, '07:18:07', '07:55', '10_OLMOS', 37, 'LP1912'),
    @(66, '07:18:07', '08:00', '23_HERNANDEZ', 42, 'LP1912'),
    @(67, '07:18:07', '08:04', '11_ETCHEVERRY', 46, 'LP1912'),
    @(68, '06:53:31', '08:05', '23_HERNANDEZ', 72, 'LP1912'),
    @(69, '07:18:07', '08:12', '15_ABASTO', 54, 'LP1912'),
    @(70, '07:18:07', '08:21', '26_HERNANDEZ', 63, 'LP1912'),
    @(71, '06:53:31', '08:22', '16_P MOR-SANTA ANA', 89, 'LP1912'),
    @(72, '07:18:07', '08:23', '16_P MOR-SANTA ANA', 65, 'LP1912'),
    @(73, '07:18:07', '08:23', '215B_EL PATO', 65, 'LP1912'),
    @(74, '07:18:07', '08:27', '84_COLONIA URQUIZA-ESC 49', 69, 'LP1912'),
    @(75, '07:18:07', '08:42', '81_EL PELIGRO', 84, 'LP1912'),
    @(76, '07:18:07', '08:44', '14_ABASTO', 86, 'LP1912'),
    @(77, '07:18:07', '08:54', '17_ROMERO', 96, 'LP1912'),
    @(78, '07:18:07', '09:02', '215A_EL PATO', 104, 'LP1912'),
    @(79, '07:18:07', '09:11', '16_P MOR-SANTA ANA', 113, 'LP1912'),
    @(80, '07:18:07', '09:17', '27_EL RETIRO', 119, 'LP1912')
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:18:07'
$ws.Cells.Item(3, 1).Value = 'Total filas: 10'

$rows = @(
    @(6, '04:37:19', '04:46', '215A_EL PATO', 9, 'LP1912'),
    @(7, '05:20:30', '05:34', '215B_EL PATO', 14, 'LP1912'),
    @(8, '04:03:00', '05:35', '215B_EL PATO', 92, 'LP1912'),
    @(9, '05:54:55', '06:11', '215A_EL PATO', 17, 'LP1912'),
    @(10, '06:24:16', '06:46', '215C_EL PATO', 22, 'LP1912'),
    @(11, '06:53:31', '07:11', '215A_EL PATO', 18, 'LP1912'),
    @(12, '06:53:31', '07:51', '215D_EL PATO', 58, 'LP1912'),
    @(13, '07:18:07', '07:52', '215D_EL PATO', 34, 'LP1912'),
    @(14, '07:18:07', '08:23', '215B_EL PATO', 65, 'LP1912'),
    @(15, '07:18:07', '09:02', '215A_EL PATO', 104, 'LP1912')
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:18:07'
$ws.Cells.Item(3, 1).Value = 'Total filas: 13'

$rows = @(
    @(6, '05:20:30', '05:44', '215A_LA PLATA', 24, 'L6173'),
    @(7, '05:54:55', '06:09', '215A_LA PLATA', 15, 'L6173'),
    @(8, '05:20:30', '06:10', '215A_LA PLATA', 50, 'L6173'),
    @(9, '05:54:55', '06:32', '215C_LA PLATA', 38, 'L6203'),
    @(10, '06:24:16', '06:33', '215C_LA PLATA', 9, 'L6203'),
    @(11, '05:54:55', '06:59', '215B_LP-P MOR-1 Y 57', 65, 'L6173'),
    @(12, '06:53:31', '07:00', '215B_LP-P MOR-1 Y 57', 7, 'L6173'),
    @(13, '05:54:55', '07:34', '215A_LA PLATA', 100, 'L6173'),
    @(14, '07:18:07', '07:35', '215A_LA PLATA', 17, 'L6173'),
    @(15, '06:53:31', '08:07', '215C_LA PLATA', 74, 'L6203'),
    @(16, '07:18:07', '08:14', '215C_LA PLATA', 56, 'L6203'),
    @(17, '07:18:07', '08:35', '215A_LA PLATA', 77, 'L6173'),
    @(18, '07:18:07', '09:09', '215D_LA PLATA', 111, 'L6203')
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
